$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# instance 1 (row 2)
$ws.Range("B2").Value = -71.65924151606299
$ws.Range("C2").Value = [double]"0.0005680745486136919"
$ws.Range("D2").Value = 28.987442266

# instance 2 (row 3)
$ws.Range("B3").Value = -69.45256735952894
$ws.Range("C3").Value = 0.0
$ws.Range("D3").Value = 77.383605264

# instance 3 (row 4)
$ws.Range("B4").Value = -71.24685956643927
$ws.Range("C4").Value = 0.05095443095246155
$ws.Range("D4").Value = 9.853651277

# instance 4 (row 5)
$ws.Range("B5").Value = -71.65440618645333
$ws.Range("C5").Value = 0.0942196409173878
$ws.Range("D5").Value = 6.576298939

# instance 5 (row 6)
$ws.Range("B6").Value = -70.22823688201865
$ws.Range("C6").Value = 0.05353691840979169
$ws.Range("D6").Value = 9.01459029

# instance 6 (row 7)
$ws.Range("B7").Value = -71.84370853874832
$ws.Range("C7").Value = 0.0
$ws.Range("D7").Value = 183.3662503

# instance 7 (row 8)
$ws.Range("B8").Value = -67.12860480761344
$ws.Range("C8").Value = [double]"5.50528395181048e-7"
$ws.Range("D8").Value = 19.50271207

# instance 8 (row 9)
$ws.Range("B9").Value = -71.56787798834289
$ws.Range("C9").Value = 0.0958882217621773
$ws.Range("D9").Value = 5.997926433

# instance 9 (row 10)
$ws.Range("B10").Value = -70.884938041055
$ws.Range("C10").Value = 0.05531266370132536
$ws.Range("D10").Value = 8.302790125

# instance 10 (row 11)
$ws.Range("B11").Value = -68.19423200083996
$ws.Range("C11").Value = 0.0
$ws.Range("D11").Value = 45.845571389
